$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.484.70"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.163.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.19"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.27%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0857"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.483.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.815"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.157.93"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.524.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.88"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.74"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.46"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0231"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.42"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.528.37"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.367.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.18%  "
